$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.971.04'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '3.018.69'
$ws.Range('E3').Value = '  -2.49%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '557.78'
$ws.Range('E5').Value = '  +0.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.07'
$ws.Range('E6').Value = '  -3.51%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.560'
$ws.Range('E8').Value = '  -1.97%  '
$ws.Range('D9').Value = '3.023.44'
$ws.Range('E9').Value = '  -2.62%  '
$ws.Range('E10').Value = '  +0.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.41'
$ws.Range('E11').Value = '  -3.97%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.366'
$ws.Range('E12').Value = '  -1.13%  '
$ws.Range('D13').Value = '3.542.66'
$ws.Range('E13').Value = '  -3.02%  '
$ws.Range('E14').Value = '  -2.73%  '
$ws.Range('D15').Value = '63.036.59'
$ws.Range('E15').Value = '  +0.07%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '24.02'
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000150'
$ws.Range('E17').Value = '  +0.45%  '
$ws.Range('D18').Value = '3.015.21'
$ws.Range('E18').Value = '  -3.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '395.93'
$ws.Range('E19').Value = '  +0.78%  '
$ws.Range('E20').Value = '  +0.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.93'
$ws.Range('E21').Value = '  -2.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.64'
$ws.Range('E22').Value = '  -3.98%  '
$ws.Range('E23').Value = '  +0.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.10'
$ws.Range('E24').Value = '  -3.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.466'
$ws.Range('E25').Value = '  -0.58%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.188'
$ws.Range('E26').Value = '  -4.85%  '
$ws.Range('D27').Value = '0.0₃0977'
$ws.Range('E27').Value = '  -0.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.68'
$ws.Range('E28').Value = '  +2.89%  '
$ws.Range('E29').Value = '  -0.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.75'
$ws.Range('E31').Value = '  -0.43%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.39'
$ws.Range('E32').Value = '  -0.98%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '159.07'
$ws.Range('E33').Value = '  +3.57%  '
$ws.Range('E34').Value = '  +0.93%  '
$ws.Range('E35').Value = '  +2.68%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.03'
$ws.Range('E36').Value = '  -0.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.30'
$ws.Range('E37').Value = '  +1.79%  '
$ws.Range('D38').Value = '2.531.50'
$ws.Range('E38').Value = '  -5.16%  '
$ws.Range('E39').Value = '  -2.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '22.69'
$ws.Range('E40').Value = '  -0.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.92'
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '37.43'
$ws.Range('E42').Value = '  -2.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.668'
$ws.Range('E43').Value = '  -2.82%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0598'
$ws.Range('E44').Value = '  +0.66%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0248'
$ws.Range('E45').Value = '  -0.55%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.997'
$ws.Range('E46').Value = '  -0.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.04'
$ws.Range('E47').Value = '  -4.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '19.95'
$ws.Range('E48').Value = '  -1.89%  '
$ws.Range('B49').Value = 'WhiteBITCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.50'
$ws.Range('E49').Value = '  +0.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0944'
$ws.Range('E50').Value = '  -1.65%  '
$ws.Range('B51').Value = 'Bittensor'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '265.06'
$ws.Range('E51').Value = '  -4.43%  '
